$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.629.16"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").Value = "2.473.89"
$ws.Range("E3").Value = "  -0.23%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'318.03"
$ws.Range("E5").Value = "  +1.45%  "
$ws.Range("D6").Value = "'92.42"
$ws.Range("E6").Value = "  -0.50%  "
$ws.Range("D7").Value = "'0.553"
$ws.Range("E7").Value = "  +1.53%  "
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("D9").Value = "'0.516"
$ws.Range("E9").Value = "  +1.76%  "
$ws.Range("D10").Value = "'0.0866"
$ws.Range("E10").Value = "  +10.52%  "
$ws.Range("D11").Value = "'33.14"
$ws.Range("E11").Value = "  +1.24%  "
$ws.Range("E12").Value = "  +0.87%  "
$ws.Range("D13").Value = "2.854.66"
$ws.Range("E13").Value = "  -0.23%  "
$ws.Range("E14").Value = "  +0.72%  "
$ws.Range("D15").Value = "'15.65"
$ws.Range("E15").Value = "  -3.15%  "
$ws.Range("D16").Value = "2.479.39"
$ws.Range("E16").Value = "  -2.64%  "
$ws.Range("D17").Value = "'0.791"
$ws.Range("E17").Value = "  +3.06%  "
$ws.Range("D18").Value = "41.593.88"
$ws.Range("E18").Value = "  +0.18%  "
$ws.Range("E19").Value = "  +1.42%  "
$ws.Range("D20").Value = "0.0₃0953"
$ws.Range("E20").Value = "  +1.52%  "
$ws.Range("D21").Value = "'71.40"
$ws.Range("E21").Value = "  -0.46%  "
$ws.Range("D22").Value = "'11.37"
$ws.Range("E22").Value = "  +0.79%  "
$ws.Range("D23").Value = "'241.01"
$ws.Range("E23").Value = "  +1.80%  "
$ws.Range("D24").Value = "'2.75"
$ws.Range("E24").Value = "  +1.97%  "
$ws.Range("E25").Value = "  +1.78%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").Value = "'24.78"
$ws.Range("E27").Value = "  -0.57%  "
$ws.Range("E28").Value = "  +3.79%  "
$ws.Range("D29").Value = "'9.88"
$ws.Range("E29").Value = "  +2.51%  "
$ws.Range("D30").Value = "'36.35"
$ws.Range("E30").Value = "  +1.10%  "
$ws.Range("D31").Value = "'159.78"
$ws.Range("E31").Value = "  +1.01%  "
$ws.Range("D32").Value = "'5.53"
$ws.Range("E33").Value = "  +0.05%  "
$ws.Range("D34").Value = "'0.0773"
$ws.Range("E34").Value = "  +2.20%  "
$ws.Range("E35").Value = "  +0.29%  "
$ws.Range("D36").Value = "'17.31"
$ws.Range("E36").Value = "  +0.32%  "
$ws.Range("E37").Value = "  -0.20%  "
$ws.Range("D38").Value = "'1.84"
$ws.Range("E38").Value = "  +0.39%  "
$ws.Range("E39").Value = "  +1.47%  "
$ws.Range("E40").Value = "  -2.23%  "
$ws.Range("E41").Value = "  -3.08%  "
$ws.Range("D42").Value = "'2.48"
$ws.Range("E42").Value = "  +1.72%  "
$ws.Range("D43").Value = "1.988.76"
$ws.Range("E43").Value = "  +0.46%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'19.19"
$ws.Range("E44").Value = "  -1.17%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "'0.0286"
$ws.Range("E45").Value = "  +1.10%  "
$ws.Range("D46").Value = "'3.00"
$ws.Range("E46").Value = "  +1.65%  "
$ws.Range("D47").Value = "'9.21"
$ws.Range("E47").Value = "  +2.73%  "
$ws.Range("D48").Value = "2.712.50"
$ws.Range("E48").Value = "  -0.26%  "
$ws.Range("D49").Value = "'97.65"
$ws.Range("E49").Value = "  -0.02%  "
$ws.Range("D50").Value = "'67.42"
$ws.Range("E50").Value = "  -1.00%  "
$ws.Range("E51").Value = "  +1.85%  "
